$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features = "['age', 'anaemia', 'creatinine_phosphokinase', 'diabetes', 'ejection_fraction', 'high_blood_pressure', 'platelets', 'serum_creatinine', 'serum_sodium', 'sex', 'smoking']"

# Header changes
$ws.Range("E1").Value = "features"
$ws.Range("H1").Value = "results_df"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Row 2 (data row 1)
$ws.Range("E2").Value = $features
$ws.Range("H2").Value = "df0"

# Row 3 (data row 2)
$ws.Range("D3").Value = $true
$ws.Range("E3").Value = $features
$ws.Range("H3").Value = "df1"

# Row 4 (data row 3)
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = $features
$ws.Range("H4").Value = "df2"

# Row 5 (data row 4)
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = $features
$ws.Range("H5").Value = "df3"

# Delete rows 6-9
$ws.Range("A6:G9").EntireRow.Delete()
